$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMap  = $wb.Worksheets.Item("Mapping Table 0")

# Update the Date property (Metadata sheet, row 8 / column B)
$wsMeta.Range("B8").Value = "2023-09-01T15:11:28+00:00"

# Correct the Source / Target group values used in the first mapping row
# (Mapping Table 0 sheet, row 2): point them at the *code system* instead
# of the *value set*.
$wsMap.Range("A2").Value = "eclaire-study-phase-source-code-system"
$wsMap.Range("D2").Value = "eclaire-study-phase-code-system"
